$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns D (Price) and E (Volume) are treated as text so that
# numeric/percent-looking strings (e.g. "1.002", "13.00", "  -0.33%  ")
# are not auto-converted to numbers by Excel when assigned via .Value
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '30.493.20'
$ws.Range('E2').Value = '  -0.33%  '
$ws.Range('D3').Value = '1.919.48'
$ws.Range('E3').Value = '  -0.18%  '
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.21%  '
$ws.Range('D5').Value = '245.83'
$ws.Range('E5').Value = '  -0.37%  '
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  +0.05%  '
$ws.Range('D7').Value = '0.4841'
$ws.Range('E7').Value = '  +2.56%  '
$ws.Range('D8').Value = '0.2899'
$ws.Range('E8').Value = '  +0.20%  '
$ws.Range('D9').Value = '0.06704'
$ws.Range('E9').Value = '  -1.28%  '
$ws.Range('D10').Value = '111.89'
$ws.Range('E10').Value = '  +5.92%  '
$ws.Range('D11').Value = '19.21'
$ws.Range('E11').Value = '  +4.36%  '
$ws.Range('D12').Value = '1.924.53'
$ws.Range('E12').Value = '  +0.49%  '
$ws.Range('D13').Value = '0.07585'
$ws.Range('E13').Value = '  -1.49%  '
$ws.Range('D14').Value = '5.343'
$ws.Range('E14').Value = '  +1.14%  '
$ws.Range('D15').Value = '0.6716'
$ws.Range('E15').Value = '  -0.20%  '
$ws.Range('D16').Value = '296.29'
$ws.Range('E16').Value = '  +1.56%  '
$ws.Range('D17').Value = '30.513.84'
$ws.Range('E17').Value = '  -0.26%  '
$ws.Range('B18').Value = 'Avalanche'
$ws.Range('C18').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D18').Value = '12.99'
$ws.Range('E18').Value = '  +0.47%  '
$ws.Range('B19').Value = 'Dai'
$ws.Range('C19').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D19').Value = '1.002'
$ws.Range('E19').Value = '  +0.21%  '
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').Value = '0.000007546'
$ws.Range('E20').Value = '  -0.88%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').Value = '5.528'
$ws.Range('E21').Value = '  +1.17%  '
$ws.Range('B22').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C22').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D22').Value = '2.172.73'
$ws.Range('E22').Value = '  +0.61%  '
$ws.Range('D23').Value = '1.002'
$ws.Range('E23').Value = '  +0.20%  '
$ws.Range('D24').Value = '6.500'
$ws.Range('E24').Value = '  +2.62%  '
$ws.Range('D25').Value = '9.449'
$ws.Range('E25').Value = '  +0.32%  '
$ws.Range('D26').Value = '164.21'
$ws.Range('E26').Value = '  -2.40%  '
$ws.Range('D27').Value = '20.47'
$ws.Range('E27').Value = '  -2.41%  '
$ws.Range('D28').Value = '2.111'
$ws.Range('E28').Value = '  -0.29%  '
$ws.Range('D29').Value = '0.1073'
$ws.Range('E29').Value = '  -0.43%  '
$ws.Range('D30').Value = '1.437'
$ws.Range('E30').Value = '  +3.28%  '
$ws.Range('D31').Value = '4.145'
$ws.Range('E31').Value = '  -0.93%  '
$ws.Range('D32').Value = '4.046'
$ws.Range('E32').Value = '  -2.17%  '
$ws.Range('D33').Value = '0.05013'
$ws.Range('E33').Value = '  -0.84%  '
$ws.Range('D34').Value = '0.7388'
$ws.Range('E34').Value = '  -0.56%  '
$ws.Range('D35').Value = '1.140'
$ws.Range('D37').Value = '2.718'
$ws.Range('E37').Value = '  -1.12%  '
$ws.Range('D38').Value = '0.02017'
$ws.Range('E38').Value = '  -3.19%  '
$ws.Range('D39').Value = '2.691'
$ws.Range('E39').Value = '  +0.00%  '
$ws.Range('D40').Value = '110.29'
$ws.Range('E40').Value = '  -1.04%  '
$ws.Range('D41').Value = '2.017'
$ws.Range('E41').Value = '  -1.77%  '
$ws.Range('D42').Value = '0.4430'
$ws.Range('E42').Value = '  +0.60%  '
$ws.Range('D43').Value = '0.8653'
$ws.Range('E43').Value = '  -1.86%  '
$ws.Range('D44').Value = '5.840'
$ws.Range('E44').Value = '  -0.80%  '
$ws.Range('D45').Value = '70.16'
$ws.Range('E45').Value = '  +4.53%  '
$ws.Range('D46').Value = '1.000'
$ws.Range('E46').Value = '  +0.05%  '
$ws.Range('E47').Value = '  -0.09%  '
$ws.Range('D48').Value = '48.76'
$ws.Range('E48').Value = '  +1.39%  '
$ws.Range('D49').Value = '9.169'
$ws.Range('E49').Value = '  -1.87%  '
$ws.Range('E50').Value = '  -1.15%  '
$ws.Range('D51').Value = '0.2504'
$ws.Range('E51').Value = '  +2.16%  '
